# Hortaliza, Agrícola del Norte S.A. de Arica - Perejil
# Weekly price update: insert a new weekly record as row 62, shifting the
# existing rows 62-64 down to 63-65 (dimension grows from R64 to R65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 62.
$ws.Rows("62:62").Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value  = 1
$ws.Cells.Item(62, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(62, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(62, 4).Value  = 45267
$ws.Cells.Item(62, 5).Value  = 15
$ws.Cells.Item(62, 6).Value  = 100112044
$ws.Cells.Item(62, 7).Value  = "Perejil"
$ws.Cells.Item(62, 8).Value  = "Sin especificar"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 300
$ws.Cells.Item(62, 11).Value = 1200
$ws.Cells.Item(62, 12).Value = 1500
$ws.Cells.Item(62, 13).Value = 1350
$ws.Cells.Item(62, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 675
$ws.Cells.Item(62, 17).Value = 2
$ws.Cells.Item(62, 18).Value = "Hortaliza"
